$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2454.853
$ws.Range("I98").Value = 2099.0386
$ws.Range("J98").Value = 3611.25
$ws.Range("K98").Value = 2099.0386
$ws.Range("L98").Value = 3611.25
$ws.Range("M98").Value = -601.0385999999999
$ws.Range("N98").Value = -6607.25
$ws.Range("H101").Value = 1644.7858
$ws.Range("I101").Value = 305.42856
$ws.Range("J101").Value = 2984.1428
$ws.Range("K101").Value = 916.28568
$ws.Range("L101").Value = 8952.428400000001
$ws.Range("M101").Value = 705.71432
$ws.Range("N101").Value = -12196.4284
$ws.Range("H106").Value = 1999.4762
$ws.Range("I106").Value = 1421
$ws.Range("K106").Value = 1421
$ws.Range("M106").Value = -790
$ws.Range("H112").Value = 1630.6923
$ws.Range("I112").Value = 2533.1667
$ws.Range("J112").Value = 857.1429000000001
$ws.Range("K112").Value = 7599.500100000001
$ws.Range("L112").Value = 2571.4287
$ws.Range("M112").Value = -6491.500100000001
$ws.Range("N112").Value = -4787.4287
$ws.Range("H122").Value = 2454.853
$ws.Range("I122").Value = 2099.0386
$ws.Range("J122").Value = 3611.25
$ws.Range("K122").Value = 6297.1158
$ws.Range("L122").Value = 10833.75
$ws.Range("M122").Value = -3847.1158
$ws.Range("N122").Value = -15733.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 987.6667
$ws.Range("I4").Value = 681.5
$ws.Range("J4").Value = 1600
$ws.Range("K4").Value = 681.5
$ws.Range("L4").Value = 1600
$ws.Range("M4").Value = -565.5
$ws.Range("N4").Value = -1832
$ws.Range("H14").Value = 400004.66
$ws.Range("I14").Value = 1000000
$ws.Range("J14").Value = 100007
$ws.Range("K14").Value = 1000000
$ws.Range("L14").Value = 100007
$ws.Range("M14").Value = -999825
$ws.Range("N14").Value = -100357
$ws.Range("H32").Value = 5089.73
$ws.Range("I32").Value = 4938.5
$ws.Range("J32").Value = 12500
$ws.Range("K32").Value = 4938.5
$ws.Range("L32").Value = 12500
$ws.Range("M32").Value = -4651.5
$ws.Range("N32").Value = -13074
$ws.Range("H113").Value = 33688.8
$ws.Range("J113").Value = 33688.8
$ws.Range("L113").Value = 33688.8
$ws.Range("N113").Value = -42366.8
$ws.Range("H122").Value = 1963.6
$ws.Range("I122").Value = 1694.25
$ws.Range("J122").Value = 2681.8667
$ws.Range("K122").Value = 5082.75
$ws.Range("L122").Value = 8045.6001
$ws.Range("M122").Value = -2632.75
$ws.Range("N122").Value = -12945.6001
$ws.Range("H132").Value = 2408.4814
$ws.Range("I132").Value = 1911.8223
$ws.Range("K132").Value = 5735.4669
$ws.Range("M132").Value = -3205.4669

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3353.9
$ws.Range("I99").Value = 2971.8667
$ws.Range("K99").Value = 2971.8667
$ws.Range("M99").Value = -1473.8667
$ws.Range("H134").Value = 2672.0293
$ws.Range("I134").Value = 2261.8696
$ws.Range("J134").Value = 3529.6365
$ws.Range("K134").Value = 6785.6088
$ws.Range("L134").Value = 10588.9095
$ws.Range("M134").Value = -4250.6088
$ws.Range("N134").Value = -15658.9095

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1699.5
$ws.Range("J15").Value = 1699.5
$ws.Range("L15").Value = 1699.5
$ws.Range("N15").Value = -2039.5
$ws.Range("H31").Value = 1861.1414
$ws.Range("I31").Value = 1211
$ws.Range("J31").Value = 2402.926
$ws.Range("K31").Value = 1211
$ws.Range("L31").Value = 2402.926
$ws.Range("M31").Value = -916
$ws.Range("N31").Value = -2992.926
$ws.Range("H34").Value = 1861.1414
$ws.Range("I34").Value = 1211
$ws.Range("J34").Value = 2402.926
$ws.Range("K34").Value = 1211
$ws.Range("L34").Value = 2402.926
$ws.Range("M34").Value = -1009
$ws.Range("N34").Value = -2806.926

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1514.0667
$ws.Range("I47").Value = 250.375
$ws.Range("K47").Value = 751.125
$ws.Range("M47").Value = -320.125
$ws.Range("H60").Value = 7964.2856
$ws.Range("I60").Value = 150
$ws.Range("J60").Value = 12305.556
$ws.Range("K60").Value = 450
$ws.Range("L60").Value = 36916.66800000001
$ws.Range("M60").Value = -199
$ws.Range("N60").Value = -37418.66800000001
$ws.Range("H75").Value = 1860.5883
$ws.Range("I75").Value = 1000
$ws.Range("J75").Value = 2045
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 6135
$ws.Range("M75").Value = -2002
$ws.Range("N75").Value = -8131
$ws.Range("H76").Value = 3856.875
$ws.Range("J76").Value = 3979.2856
$ws.Range("L76").Value = 11937.8568
$ws.Range("N76").Value = -12703.8568
$ws.Range("H78").Value = 1860.5883
$ws.Range("I78").Value = 1000
$ws.Range("J78").Value = 2045
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 18405
$ws.Range("M78").Value = -4008
$ws.Range("N78").Value = -28389
$ws.Range("H79").Value = 3856.875
$ws.Range("J79").Value = 3979.2856
$ws.Range("L79").Value = 11937.8568
$ws.Range("N79").Value = -14589.8568
$ws.Range("H87").Value = 7592.643
$ws.Range("I87").Value = 2921.889
$ws.Range("K87").Value = 8765.667000000001
$ws.Range("M87").Value = -7517.667000000001
$ws.Range("H90").Value = 7592.643
$ws.Range("I90").Value = 2921.889
$ws.Range("K90").Value = 26297.001
$ws.Range("M90").Value = -20057.001
$ws.Range("H94").Value = 3497
$ws.Range("I94").Value = 1331
$ws.Range("J94").Value = 3806.4285
$ws.Range("K94").Value = 3993
$ws.Range("L94").Value = 11419.2855
$ws.Range("M94").Value = -3317
$ws.Range("N94").Value = -12771.2855
$ws.Range("H125").Value = 2648.1
$ws.Range("I125").Value = 1900
$ws.Range("J125").Value = 2835.125
$ws.Range("K125").Value = 5700
$ws.Range("L125").Value = 8505.375
$ws.Range("M125").Value = -780
$ws.Range("N125").Value = -18345.375
$ws.Range("H131").Value = 1224.1177
$ws.Range("I131").Value = 3080
$ws.Range("J131").Value = 1031.2987
$ws.Range("K131").Value = 9240
$ws.Range("L131").Value = 3093.8961
$ws.Range("M131").Value = -4200
$ws.Range("N131").Value = -13173.8961

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5236.5
$ws.Range("I122").Value = 4481.5454
$ws.Range("J122").Value = 6422.857
$ws.Range("K122").Value = 13444.6362
$ws.Range("L122").Value = 19268.571
$ws.Range("M122").Value = -10994.6362
$ws.Range("N122").Value = -24168.571
$ws.Range("H132").Value = 4065.5227
$ws.Range("I132").Value = 4415.4
$ws.Range("J132").Value = 3605.158
$ws.Range("K132").Value = 13246.2
$ws.Range("L132").Value = 10815.474
$ws.Range("M132").Value = -10716.2
$ws.Range("N132").Value = -15875.474

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5038.3076
$ws.Range("I40").Value = 6500
$ws.Range("K40").Value = 6500
$ws.Range("M40").Value = -6364
$ws.Range("H100").Value = 3160
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3160
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 3160
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -4242
$ws.Range("H106").Value = 26330.908
$ws.Range("J106").Value = 26330.908
$ws.Range("L106").Value = 26330.908
$ws.Range("N106").Value = -28854.908
$ws.Range("H122").Value = 2857.1
$ws.Range("I122").Value = 2346.1538
$ws.Range("K122").Value = 7038.4614
$ws.Range("M122").Value = -4588.4614
$ws.Range("H132").Value = 2543.1143
$ws.Range("I132").Value = 1633.2667
$ws.Range("J132").Value = 3225.5
$ws.Range("K132").Value = 4899.800099999999
$ws.Range("L132").Value = 9676.5
$ws.Range("M132").Value = -2369.800099999999
$ws.Range("N132").Value = -14736.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 27642
$ws.Range("J98").Value = 27642
$ws.Range("L98").Value = 27642
$ws.Range("N98").Value = -33632
$ws.Range("H122").Value = 3279.8
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050
$ws.Range("H132").Value = 6850.1333
$ws.Range("I132").Value = 4147.7334
$ws.Range("K132").Value = 12443.2002
$ws.Range("M132").Value = -9913.200199999999
